$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 56
$srcRow = 55

# Column A holds a date-formatted-looking literal ("2026/01/05") that must stay
# stored as literal text (not get auto-parsed into a date serial number) -
# route it through a TEXT() formula in a scratch cell, then paste the
# resulting value back in as a plain value (no formula left behind).
$scratch = $ws.Cells.Item(1000, 26)
$scratch.Formula = '=TEXT("2026/01/05","@")'
$scratch.Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4163)
$scratch.Clear()

$ws.Cells.Item($row, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($row, 3).Value = 1135

# Copy the formatting (style) of the previous row's cells onto the new row,
# matching s="1" without minting any new/orphan cell-style records.
$ws.Range($ws.Cells.Item($srcRow, 1), $ws.Cells.Item($srcRow, 3)).Copy()
$ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 3)).PasteSpecial(-4122)
